$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$val1 = $ws.Range("A1").Text
$val2 = $ws.Range("A2").Text
$val3 = $ws.Range("A3").Text
$val4 = $ws.Range("A4").Text

for ($i = 0; $i -lt 3; $i++) {
    $baseRow = 5 + ($i * 4)
    $ws.Range("A$baseRow").Value = $val1
    $ws.Range("A$($baseRow+1)").Value = $val2
    $ws.Range("A$($baseRow+2)").Value = $val3
    $ws.Range("A$($baseRow+3)").Value = $val4
}

$ws.Range("E7").Select()
